$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Update the author/name placeholder text.
$s.Shapes.Item(2).TextFrame.TextRange.Text = " | WATKIDOG JOACHIM | MUGISHA RONALD | LOCHA DERRICK | ACAA GLADYS OBOL "

# 2. Fill in the Project Requirements content (R1-R4).
$reqText = "R1.The system must be able to detect temperature, light, moisture level, humidity within a greenhouse.`rR2. The system must store the data collected in a database.`rR3. The system must use the data collected in the database to automatically control the conditions within a greenhouse.`rR4. The system must be able to analyze the data collected."
$s.Shapes.Item(8).TextFrame.TextRange.Text = $reqText

# 3. Update the "Target Users" SmartArt labels to "User category xx".
$sa = $s.Shapes.Item(10).SmartArt
$allNodes = $sa.AllNodes
$allNodes.Item(1).TextFrame2.TextRange.Text = "User category xx"
$allNodes.Item(3).TextFrame2.TextRange.Text = "User category xx"
$allNodes.Item(5).TextFrame2.TextRange.Text = "User category xx"

# 4. Reposition/resize the logo picture (top-right image).
$pic = $s.Shapes.Item(23)
$pic.Left = 2266.8387401574805
$pic.Width = 1189.1613385826772

# 5. Nudge the title placeholder box to match the recalculated layout.
$title = $s.Shapes.Item(1)
$title.Left = 1.7419685039370079
$title.Top = 46.676062992125985
$title.Width = 2376.0
$title.Height = 233.9952755905512
